$wb = $excel.ActiveWorkbook

# ==== Sheet: LP1912 ====
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 19:54:57'
$ws.Cells.Item(3, 1).Value = 'Total filas: 344'
$ws.Cells.Item(133, 1).Value = '10:36:50'
$ws.Cells.Item(133, 3).Value = '14_ABASTO'
$ws.Cells.Item(133, 4).Value = 116
$ws.Cells.Item(134, 1).Value = '11:33:52'
$ws.Cells.Item(134, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(134, 4).Value = 59
$ws.Cells.Item(144, 1).Value = '10:49:38'
$ws.Cells.Item(144, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(144, 4).Value = 119
$ws.Cells.Item(145, 1).Value = '11:33:52'
$ws.Cells.Item(145, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(145, 4).Value = 75
$ws.Cells.Item(210, 1).Value = '15:16:46'
$ws.Cells.Item(210, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(210, 4).Value = 49
$ws.Cells.Item(211, 1).Value = '14:11:28'
$ws.Cells.Item(211, 3).Value = '14_ABASTO'
$ws.Cells.Item(211, 4).Value = 114
$ws.Cells.Item(220, 1).Value = '16:12:06'
$ws.Cells.Item(220, 3).Value = '14_ABASTO'
$ws.Cells.Item(220, 4).Value = 18
$ws.Cells.Item(221, 1).Value = '15:16:46'
$ws.Cells.Item(221, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(221, 4).Value = 74
$ws.Cells.Item(294, 1).Value = '17:55:25'
$ws.Cells.Item(294, 3).Value = '14_ABASTO'
$ws.Cells.Item(294, 4).Value = 68
$ws.Cells.Item(295, 1).Value = '17:35:41'
$ws.Cells.Item(295, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(295, 4).Value = 88
$ws.Cells.Item(317, 1).Value = '17:55:25'
$ws.Cells.Item(317, 3).Value = '225_GOMEZ'
$ws.Cells.Item(317, 4).Value = 118
$ws.Cells.Item(318, 1).Value = '18:52:29'
$ws.Cells.Item(318, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(318, 4).Value = 61
$ws.Cells.Item(319, 1).Value = '19:54:57'
$ws.Cells.Item(319, 2).Value = '19:54'
$ws.Cells.Item(319, 3).Value = '225_GOMEZ'
$ws.Cells.Item(319, 4).Value = 0
$ws.Cells.Item(320, 1).Value = '19:54:57'
$ws.Cells.Item(320, 2).Value = '20:04'
$ws.Cells.Item(320, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(320, 4).Value = 10
$ws.Cells.Item(321, 1).Value = '19:35:34'
$ws.Cells.Item(321, 2).Value = '20:05'
$ws.Cells.Item(321, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(321, 4).Value = 30
$ws.Cells.Item(322, 1).Value = '18:11:09'
$ws.Cells.Item(322, 2).Value = '20:06'
$ws.Cells.Item(322, 3).Value = '215C_EL PATO'
$ws.Cells.Item(322, 4).Value = 115
$ws.Cells.Item(323, 1).Value = '18:52:29'
$ws.Cells.Item(323, 2).Value = '20:07'
$ws.Cells.Item(323, 3).Value = '215C_EL PATO'
$ws.Cells.Item(323, 4).Value = 75
$ws.Cells.Item(324, 1).Value = '18:52:29'
$ws.Cells.Item(324, 2).Value = '20:09'
$ws.Cells.Item(324, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(324, 4).Value = 77
$ws.Cells.Item(325, 1).Value = '18:44:45'
$ws.Cells.Item(325, 2).Value = '20:10'
$ws.Cells.Item(325, 4).Value = 86
$ws.Cells.Item(326, 1).Value = '18:44:45'
$ws.Cells.Item(326, 2).Value = '20:12'
$ws.Cells.Item(326, 3).Value = '14_ABASTO'
$ws.Cells.Item(326, 4).Value = 88
$ws.Cells.Item(327, 1).Value = '19:11:44'
$ws.Cells.Item(327, 2).Value = '20:13'
$ws.Cells.Item(327, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(327, 4).Value = 62
$ws.Cells.Item(328, 1).Value = '18:30:48'
$ws.Cells.Item(328, 2).Value = '20:21'
$ws.Cells.Item(328, 3).Value = '15_ABASTO'
$ws.Cells.Item(328, 4).Value = 111
$ws.Cells.Item(329, 1).Value = '18:52:29'
$ws.Cells.Item(329, 2).Value = '20:22'
$ws.Cells.Item(329, 3).Value = '15_ABASTO'
$ws.Cells.Item(329, 4).Value = 90
$ws.Cells.Item(330, 1).Value = '18:44:45'
$ws.Cells.Item(330, 2).Value = '20:30'
$ws.Cells.Item(330, 3).Value = '10_OLMOS'
$ws.Cells.Item(330, 4).Value = 106
$ws.Cells.Item(331, 1).Value = '19:47:50'
$ws.Cells.Item(331, 2).Value = '20:33'
$ws.Cells.Item(331, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(331, 4).Value = 46
$ws.Cells.Item(332, 1).Value = '19:35:34'
$ws.Cells.Item(332, 2).Value = '20:34'
$ws.Cells.Item(332, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(332, 4).Value = 59
$ws.Cells.Item(333, 1).Value = '19:11:44'
$ws.Cells.Item(333, 2).Value = '20:41'
$ws.Cells.Item(333, 4).Value = 90
$ws.Cells.Item(334, 1).Value = '18:52:29'
$ws.Cells.Item(334, 2).Value = '20:42'
$ws.Cells.Item(334, 4).Value = 110
$ws.Cells.Item(335, 1).Value = '19:35:34'
$ws.Cells.Item(335, 2).Value = '20:43'
$ws.Cells.Item(335, 3).Value = '17_ROMERO'
$ws.Cells.Item(335, 4).Value = 68
$ws.Cells.Item(336, 1).Value = '19:47:50'
$ws.Cells.Item(336, 2).Value = '20:45'
$ws.Cells.Item(336, 3).Value = '17_ROMERO'
$ws.Cells.Item(336, 4).Value = 58
$ws.Cells.Item(337, 1).Value = '18:52:29'
$ws.Cells.Item(337, 2).Value = '20:47'
$ws.Cells.Item(337, 3).Value = '215B_EL PATO'
$ws.Cells.Item(337, 4).Value = 115
$ws.Cells.Item(338, 1).Value = '19:35:34'
$ws.Cells.Item(338, 2).Value = '20:55'
$ws.Cells.Item(338, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(338, 4).Value = 80
$ws.Cells.Item(339, 1).Value = '19:54:57'
$ws.Cells.Item(339, 2).Value = '20:55'
$ws.Cells.Item(339, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(339, 4).Value = 61
$ws.Cells.Item(340, 1).Value = '19:11:44'
$ws.Cells.Item(340, 2).Value = '20:56'
$ws.Cells.Item(340, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(340, 4).Value = 105
$ws.Cells.Item(341, 1).Value = '19:11:44'
$ws.Cells.Item(341, 2).Value = '21:06'
$ws.Cells.Item(341, 3).Value = '10_OLMOS'
$ws.Cells.Item(341, 4).Value = 115
$ws.Cells.Item(342, 2).Value = '21:09'
$ws.Cells.Item(342, 3).Value = '15_ABASTO'
$ws.Cells.Item(342, 4).Value = 82
$ws.Cells.Item(343, 2).Value = '21:10'
$ws.Cells.Item(343, 3).Value = '15_ABASTO'
$ws.Cells.Item(343, 4).Value = 95
$ws.Cells.Item(344, 1).Value = '19:35:34'
$ws.Cells.Item(344, 2).Value = '21:28'
$ws.Cells.Item(344, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(344, 4).Value = 113
$ws.Cells.Item(345, 1).Value = '19:47:50'
$ws.Cells.Item(345, 2).Value = '21:33'
$ws.Cells.Item(345, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(345, 4).Value = 106
$ws.Cells.Item(345, 5).Value = 'LP1912'
$ws.Cells.Item(346, 1).Value = '19:54:57'
$ws.Cells.Item(346, 2).Value = '21:33'
$ws.Cells.Item(346, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(346, 4).Value = 99
$ws.Cells.Item(346, 5).Value = 'LP1912'
$ws.Cells.Item(347, 1).Value = '19:35:34'
$ws.Cells.Item(347, 2).Value = '21:34'
$ws.Cells.Item(347, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(347, 4).Value = 119
$ws.Cells.Item(347, 5).Value = 'LP1912'
$ws.Cells.Item(348, 1).Value = '19:54:57'
$ws.Cells.Item(348, 2).Value = '21:44'
$ws.Cells.Item(348, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(348, 4).Value = 110
$ws.Cells.Item(348, 5).Value = 'LP1912'
$ws.Cells.Item(349, 1).Value = '19:47:50'
$ws.Cells.Item(349, 2).Value = '21:45'
$ws.Cells.Item(349, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(349, 4).Value = 118
$ws.Cells.Item(349, 5).Value = 'LP1912'

# ==== Sheet: LP1912-215 ====
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 19:54:57'

# ==== Sheet: 6203-6173 ====
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 19:54:57'
$ws.Cells.Item(3, 1).Value = 'Total filas: 49'
$ws.Cells.Item(53, 1).Value = '19:54:57'
$ws.Cells.Item(53, 2).Value = '21:29'
$ws.Cells.Item(53, 4).Value = 95
$ws.Cells.Item(54, 1).Value = '19:35:34'
$ws.Cells.Item(54, 2).Value = '21:30'
$ws.Cells.Item(54, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(54, 4).Value = 115
$ws.Cells.Item(54, 5).Value = 'L6203'
